# Apply updated FFXIV Leve profit figures (scheduled price-refresh run)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 643.2941
$ws.Range("I6").Value = 302.4
$ws.Range("J6").Value = 3200
$ws.Range("K6").Value = 907.1999999999999
$ws.Range("L6").Value = 9600
$ws.Range("M6").Value = -795.1999999999999
$ws.Range("N6").Value = -9824

# Row 15
$ws.Range("H15").Value = 2737.78
$ws.Range("I15").Value = 2737.78
$ws.Range("K15").Value = 8213.34
$ws.Range("M15").Value = -8044.34

# Row 21
$ws.Range("H21").Value = 12103.4
$ws.Range("I21").Value = 12103.4
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 12103.4
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -11635.4
$ws.Range("N21").ClearContents()

# Row 23
$ws.Range("H23").Value = 12103.4
$ws.Range("I23").Value = 12103.4
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 12103.4
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -11869.4
$ws.Range("N23").ClearContents()

# Row 40
$ws.Range("H40").Value = 9762.333000000001
$ws.Range("I40").Value = 3980
$ws.Range("J40").Value = 10102.471
$ws.Range("K40").Value = 3980
$ws.Range("L40").Value = 10102.471
$ws.Range("M40").Value = -3805
$ws.Range("N40").Value = -10452.471

# Row 107
$ws.Range("H107").Value = 482.30768
$ws.Range("I107").Value = 504.33334
$ws.Range("J107").Value = 218
$ws.Range("K107").Value = 504.33334
$ws.Range("L107").Value = 218
$ws.Range("M107").Value = 1415.66666
$ws.Range("N107").Value = -4058

# Row 127
$ws.Range("H127").Value = 500.4
$ws.Range("I127").Value = 500.4
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1501.2
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 3458.8
$ws.Range("N127").ClearContents()

# Row 132
$ws.Range("H132").Value = 1432.5667
$ws.Range("I132").Value = 1464.862
$ws.Range("J132").Value = 496
$ws.Range("K132").Value = 4394.586
$ws.Range("L132").Value = 1488
$ws.Range("M132").Value = -1864.586
$ws.Range("N132").Value = -6548

# Row 137
$ws.Range("H137").Value = 66672092
$ws.Range("I137").Value = 142860600
$ws.Range("J137").Value = 7144.125
$ws.Range("K137").Value = 428581800
$ws.Range("L137").Value = 21432.375
$ws.Range("M137").Value = -428579250
$ws.Range("N137").Value = -26532.375

# Row 138
$ws.Range("H138").Value = 5397.2896
$ws.Range("J138").Value = 5963.815
$ws.Range("L138").Value = 17891.445
$ws.Range("N138").Value = -28171.445

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 288.33334
$ws.Range("I5").Value = 228.14285
$ws.Range("J5").Value = 499
$ws.Range("K5").Value = 228.14285
$ws.Range("L5").Value = 499
$ws.Range("M5").Value = -116.14285
$ws.Range("N5").Value = -723

# Row 132
$ws.Range("H132").Value = 4254.522
$ws.Range("I132").Value = 2945.4707
$ws.Range("J132").Value = 7963.5
$ws.Range("K132").Value = 8836.4121
$ws.Range("L132").Value = 23890.5
$ws.Range("M132").Value = -6306.4121
$ws.Range("N132").Value = -28950.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 288.33334
$ws.Range("I4").Value = 228.14285
$ws.Range("J4").Value = 499
$ws.Range("K4").Value = 228.14285
$ws.Range("L4").Value = 499
$ws.Range("M4").Value = -113.14285
$ws.Range("N4").Value = -729

# Row 135
$ws.Range("H135").Value = 59998.4
$ws.Range("J135").Value = 59998.4
$ws.Range("L135").Value = 59998.4
$ws.Range("N135").Value = -70138.39999999999

# Row 138
$ws.Range("H138").Value = 65096.8
$ws.Range("J138").Value = 65096.8
$ws.Range("L138").Value = 65096.8
$ws.Range("N138").Value = -75376.8

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 31065.309
$ws.Range("I31").Value = 3875.4614
$ws.Range("J31").Value = 85445
$ws.Range("K31").Value = 3875.4614
$ws.Range("L31").Value = 85445
$ws.Range("M31").Value = -3580.4614
$ws.Range("N31").Value = -86035

# Row 34
$ws.Range("H34").Value = 31065.309
$ws.Range("I34").Value = 3875.4614
$ws.Range("J34").Value = 85445
$ws.Range("K34").Value = 3875.4614
$ws.Range("L34").Value = 85445
$ws.Range("M34").Value = -3673.4614
$ws.Range("N34").Value = -85849

# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 5686189.5
$ws.Range("I68").Value = 11365040
$ws.Range("J68").Value = 7339.091
$ws.Range("K68").Value = 34095120
$ws.Range("L68").Value = 22017.273
$ws.Range("M68").Value = -34094309
$ws.Range("N68").Value = -23639.273

# Row 71
$ws.Range("H71").Value = 5686189.5
$ws.Range("I71").Value = 11365040
$ws.Range("J71").Value = 7339.091
$ws.Range("K71").Value = 102285360
$ws.Range("L71").Value = 66051.819
$ws.Range("M71").Value = -102281304
$ws.Range("N71").Value = -74163.819

# Row 107
$ws.Range("H107").Value = 71430900
$ws.Range("I107").Value = 1750.875
$ws.Range("J107").Value = 166669760
$ws.Range("K107").Value = 5252.625
$ws.Range("L107").Value = 500009280
$ws.Range("M107").Value = -3332.625
$ws.Range("N107").Value = -500013120

# Row 112
$ws.Range("H112").Value = 1750
$ws.Range("J112").Value = 1750
$ws.Range("L112").Value = 5250
$ws.Range("N112").Value = -7466

# Row 132
$ws.Range("H132").Value = 4305.7417
$ws.Range("I132").Value = 5051.3335
$ws.Range("J132").Value = 4000.7273
$ws.Range("K132").Value = 45462.0015
$ws.Range("L132").Value = 36006.5457
$ws.Range("M132").Value = -42932.0015
$ws.Range("N132").Value = -41066.5457

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 14380
$ws.Range("I57").Value = 3966.6667
$ws.Range("K57").Value = 3966.6667
$ws.Range("M57").Value = -3146.6667

# Row 102
$ws.Range("H102").Value = 4535.6
$ws.Range("I102").Value = 3543.923
$ws.Range("K102").Value = 3543.923
$ws.Range("M102").Value = -1921.923

$ws = $wb.Worksheets.Item("LTW")
# Row 34
$ws.Range("H34").Value = 19674.666
$ws.Range("J34").Value = 19674.666
$ws.Range("L34").Value = 19674.666
$ws.Range("N34").Value = -20018.666

# Row 40
$ws.Range("H40").Value = 6529.5
$ws.Range("I40").Value = 3764.1428
$ws.Range("J40").Value = 10401
$ws.Range("K40").Value = 3764.1428
$ws.Range("L40").Value = 10401
$ws.Range("M40").Value = -3628.1428
$ws.Range("N40").Value = -10673

# Row 61
$ws.Range("H61").Value = 5590.778
$ws.Range("I61").Value = 5590.778
$ws.Range("K61").Value = 5590.778
$ws.Range("M61").Value = -5388.778

# Row 68
$ws.Range("H68").Value = 4954.1816
$ws.Range("I68").Value = 1999.125
$ws.Range("J68").Value = 12834.333
$ws.Range("K68").Value = 1999.125
$ws.Range("L68").Value = 12834.333
$ws.Range("M68").Value = -1250.125
$ws.Range("N68").Value = -14332.333

# Row 71
$ws.Range("H71").Value = 4954.1816
$ws.Range("I71").Value = 1999.125
$ws.Range("J71").Value = 12834.333
$ws.Range("K71").Value = 9995.625
$ws.Range("L71").Value = 64171.665
$ws.Range("M71").Value = -6251.625
$ws.Range("N71").Value = -71659.66500000001

# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# Row 113
$ws.Range("H113").Value = 5590.778
$ws.Range("I113").Value = 5590.778
$ws.Range("K113").Value = 5590.778
$ws.Range("M113").Value = -3420.778

# Row 115
$ws.Range("H115").Value = 70000
$ws.Range("J115").Value = 70000
$ws.Range("L115").Value = 70000
$ws.Range("N115").Value = -72350

$ws = $wb.Worksheets.Item("WVR")
# Row 44
$ws.Range("H44").Value = 21999.2
$ws.Range("J44").Value = 21999.2
$ws.Range("L44").Value = 21999.2
$ws.Range("N44").Value = -23107.2
